# "US 3.3 commit files" — advance the IT (Initial Time) policy-schedule
# input from 2017 to 2020, and leave the workbook focused on the IT sheet
# (active tab / selection), matching the author's saved UI state.

$wb = $excel.ActiveWorkbook

$itSheet = $wb.Worksheets.Item("IT")

# Update the Initial Time value (IT!B2) from 2017 to 2020.
$itSheet.Range("B2").Value = 2020

# Selecting the cell switches to/activates its sheet, mirroring the
# commit's workbook.xml (activeTab moves to the IT sheet) and sheet2.xml
# (tabSelected + active cell B3) changes.
$itSheet.Range("B3").Select()
